# Add new drug/diagnosis/disease-code rows (227-252) uploaded in this commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(227, 1).Value = "PRIMOLUT- N 5  MG ( O )"
$ws.Cells.Item(227, 1).WrapText = $true
$ws.Cells.Item(227, 1).VerticalAlignment = -4108
$ws.Cells.Item(227, 4).Value = "E288"
$ws.Cells.Item(228, 1).Value = "NAPROXEN 250 MG"
$ws.Cells.Item(228, 1).WrapText = $true
$ws.Cells.Item(228, 1).VerticalAlignment = -4108
$ws.Cells.Item(228, 4).Value = "M543"
$ws.Cells.Item(229, 1).Value = "NORFLOXACIN 400 MG"
$ws.Cells.Item(229, 1).WrapText = $true
$ws.Cells.Item(229, 1).VerticalAlignment = -4108
$ws.Cells.Item(229, 4).Value = "N390"
$ws.Cells.Item(230, 1).Value = "GENTAMICIN CREAM (per G)"
$ws.Cells.Item(230, 1).WrapText = $true
$ws.Cells.Item(230, 1).VerticalAlignment = -4108
$ws.Cells.Item(230, 4).Value = "N390"
$ws.Cells.Item(231, 1).Value = "NORFLOXACIN 400 MG"
$ws.Cells.Item(231, 1).WrapText = $true
$ws.Cells.Item(231, 1).VerticalAlignment = -4108
$ws.Cells.Item(231, 4).Value = "A099"
$ws.Cells.Item(232, 1).Value = "PARACETAMOL 500  MG"
$ws.Cells.Item(232, 1).WrapText = $true
$ws.Cells.Item(232, 1).VerticalAlignment = -4108
$ws.Cells.Item(232, 4).Value = "A099"
$ws.Cells.Item(233, 1).Value = "TRAMAL 50  MG  CAP"
$ws.Cells.Item(233, 1).WrapText = $true
$ws.Cells.Item(233, 1).VerticalAlignment = -4108
$ws.Cells.Item(233, 4).Value = "M791"
$ws.Cells.Item(234, 1).Value = "CANESTEN VAG SUPPO 100  MG"
$ws.Cells.Item(234, 1).WrapText = $true
$ws.Cells.Item(234, 1).VerticalAlignment = -4108
$ws.Cells.Item(234, 4).Value = "B373"
$ws.Cells.Item(235, 1).Value = "TRANXENE 5 MG"
$ws.Cells.Item(235, 1).WrapText = $true
$ws.Cells.Item(235, 1).VerticalAlignment = -4108
$ws.Cells.Item(235, 4).Value = "G470"
$ws.Cells.Item(236, 1).Value = "DEXTRO 15 MG  "
$ws.Cells.Item(236, 1).WrapText = $true
$ws.Cells.Item(236, 1).VerticalAlignment = -4108
$ws.Cells.Item(236, 4).Value = "J449"
$ws.Cells.Item(237, 1).Value = "DURALYN CR 200 mg"
$ws.Cells.Item(237, 1).WrapText = $true
$ws.Cells.Item(237, 1).VerticalAlignment = -4108
$ws.Cells.Item(237, 4).Value = "J449"
$ws.Cells.Item(238, 1).Value = "PREDNISOLONE 5 MG  "
$ws.Cells.Item(238, 1).WrapText = $true
$ws.Cells.Item(238, 1).VerticalAlignment = -4108
$ws.Cells.Item(238, 4).Value = "J449"
$ws.Cells.Item(239, 1).Value = "CLARITYNE 10 MG TAB*** SA4"
$ws.Cells.Item(239, 1).WrapText = $true
$ws.Cells.Item(239, 1).VerticalAlignment = -4108
$ws.Cells.Item(239, 4).Value = "J449"
$ws.Cells.Item(240, 1).Value = "SIBELIUM  5  MG"
$ws.Cells.Item(240, 2).Value = "บรรเทาอาการบ้านหมุน"
$ws.Cells.Item(240, 4).Value = "R42"
$ws.Cells.Item(241, 1).Value = "MOTILIUM 10 MG TAB"
$ws.Cells.Item(241, 1).WrapText = $true
$ws.Cells.Item(241, 1).VerticalAlignment = -4108
$ws.Cells.Item(241, 2).Value = "คลื่นไส้"
$ws.Cells.Item(241, 4).Value = "N946"
$ws.Cells.Item(242, 1).Value = "PROCTOSEDYL  SUPPO"
$ws.Cells.Item(242, 1).WrapText = $true
$ws.Cells.Item(242, 4).Value = "K643"
$ws.Cells.Item(243, 1).Value = "DAFLON 500 MG"
$ws.Cells.Item(243, 1).WrapText = $true
$ws.Cells.Item(243, 1).VerticalAlignment = -4108
$ws.Cells.Item(243, 4).Value = "K643"
$ws.Cells.Item(244, 4).Value = "F510"
$ws.Cells.Item(244, 1).Value = "XANAX 0.25 MG***SA6"
$ws.Cells.Item(244, 1).WrapText = $true
$ws.Cells.Item(244, 1).VerticalAlignment = -4108
$ws.Cells.Item(245, 1).Value = "ATIVAN 0.5 MG*** SA4"
$ws.Cells.Item(245, 1).WrapText = $true
$ws.Cells.Item(245, 1).VerticalAlignment = -4108
$ws.Cells.Item(245, 4).Value = "F413"
$ws.Cells.Item(246, 4).Value = "F413"
$ws.Cells.Item(246, 1).Value = "ZYMRON 15 mg"
$ws.Cells.Item(247, 1).Value = "BRUFEN 400 MG "
$ws.Cells.Item(247, 1).WrapText = $true
$ws.Cells.Item(247, 1).VerticalAlignment = -4108
$ws.Cells.Item(247, 4).Value = "R51"
$ws.Cells.Item(248, 1).Value = "MESTINON  60  MG"
$ws.Cells.Item(248, 4).Value = "G700"
$ws.Cells.Item(249, 1).Value = "3TC  150  MG TAB"
$ws.Cells.Item(249, 1).WrapText = $true
$ws.Cells.Item(249, 1).VerticalAlignment = -4108
$ws.Cells.Item(249, 4).Value = "B181"
$ws.Cells.Item(250, 1).Value = "VALOSINE SR 75 mg"
$ws.Cells.Item(250, 1).WrapText = $true
$ws.Cells.Item(250, 1).VerticalAlignment = -4108
$ws.Cells.Item(251, 1).Value = "TRYPTANAL 25 MG*** LA1"
$ws.Cells.Item(251, 1).WrapText = $true
$ws.Cells.Item(251, 1).VerticalAlignment = -4108
$ws.Cells.Item(251, 4).Value = "F320"
$ws.Cells.Item(250, 4).Value = "F320"
$ws.Cells.Item(252, 4).Value = "M159"
$ws.Cells.Item(252, 1).Value = "VOLTAREN  25  MG TAB*** LA1/SA5"
$ws.Cells.Item(252, 1).WrapText = $true

# Row heights for the two wrapped long-name rows
$ws.Rows.Item(234).RowHeight = 45
$ws.Rows.Item(252).RowHeight = 45

# Restore the selection to match the author's saved session
$ws.Range("D234").Select() | Out-Null
